$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The acknowledgments table (A1:F39) is sorted alphabetically by column A.
# Row 10 is the "lemmalist-greek" dependency entry, which is being removed
# entirely (commit: "Remove lemmalist-greek"). Deleting the whole row shifts
# every row below it up by one.

# 1) Drop the two hyperlinks anchored in row 10 (home page in B10, license
#    URL in F10) before the row shifts, deleting from the highest address
#    down so earlier deletions don't renumber/shift the ones still pending.
$rowHyperlinks = @()
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$10' -or $addr -eq '$F$10') {
        $rowHyperlinks += $hl
    }
}
for ($i = $rowHyperlinks.Count - 1; $i -ge 0; $i--) {
    $rowHyperlinks[$i].Delete()
}

# 2) Delete the entire row, shifting rows 11:39 up to 10:38.
$ws.Rows.Item(10).Delete()

# 3) Re-apply the sort over the now-shrunk range so the sheet's recorded
#    sort state (A2:F38 / A2:A38) reflects the new extent.
$sortRange = $ws.Range("A2:F38")
$keyRange = $ws.Range("A2:A38")
$ws.Sort.SetRange($sortRange)
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange)
$ws.Sort.Apply()
